$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ImageName") holds values like "MEETING_2022_03_27-img1.png".
# Strip the "img" prefix and the ".png" extension, leaving
# "MEETING_2022_03_27-1".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 103 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old -replace "img", ""
        $new = $new -replace "\.png$", ""
        $cell.Value = $new
    }
}
